# Actualización desde MV -datos-
# Append the new daily "Spot" observations (27-09-2021 .. 01-10-2021) to the
# bottom of the data table on Sheet1 (columns A:E), continuing directly after
# the existing last row (188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("27-09-2021", 4.73, 3.93, 3.47, 3.37),
    @("28-09-2021", 4.78, 3.95, 3.49, 3.38),
    @("29-09-2021", 5.01, 4.05, 3.55, 3.44),
    @("30-09-2021", 5.03, 4.05, 3.51, 3.37),
    @("01-10-2021", 4.84, 4.08, 3.52, 3.37)
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
